$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers (e.g. "213.21") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers.
$textForceAddrs = @("D5", "D6", "D8", "D11", "D14", "D15", "D16", "D18", "D20", "D22", "D23", "D25", "D27", "D28", "D29", "D37", "D38", "D40", "D43", "D49", "D50", "D51")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.720.23'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.646.00'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '213.21'
$ws.Range('D6').Value = '0.531'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '23.30'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.878.58'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '1.654.05'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').Value = '0.559'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '64.73'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '27.702.36'
$ws.Range('E17').Value = '  +1.19%  '
$ws.Range('D18').Value = '231.36'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '7.63'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = '10.13'
$ws.Range('E23').Value = '  +9.23%  '
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('D25').Value = '150.15'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').Value = '0.112'
$ws.Range('E27').Value = '  -2.43%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '15.65'
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +1.18%  '
$ws.Range('D33').Value = '1.441.25'
$ws.Range('E33').Value = '  +3.17%  '
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('D37').Value = '0.569'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('D38').Value = '0.878'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('D40').Value = '0.885'
$ws.Range('E40').Value = '  +12.55%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').Value = '67.08'
$ws.Range('E43').Value = '  +4.48%  '
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('E45').Value = '  +1.74%  '
$ws.Range('D46').Value = '1.787.93'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  +6.27%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('D49').Value = '85.51'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').Value = '0.0988'
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('D51').Value = '7.75'
$ws.Range('E51').Value = '  +1.53%  '

# Reset style pointer to Normal so the cell style index matches the baseline
# (keeps the cell text-typed without leaving a stray explicit NumberFormat style).
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).Style = "Normal"
}
